$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room for the new "Meeting minutes" mini-table (rows 16-17) and for
#     the week-3 block moving from rows 17-20 down to rows 27-30 (plus a new
#     row 31 label). Inserting 10 rows at 15-24 shifts old row 17 -> row 27,
#     old rows 18/19/20 -> 28/29/30, while rows 4-14 stay untouched.
$ws.Range("A15:A24").EntireRow.Insert() | Out-Null

# --- Week2 (first) data row: tweak the existing note text.
$ws.Range("C12").Value = "10.30-> 4"

# --- Row 14 meeting-notes line shifts one "slot" to the right (A->dropped,
#     C->A, D->C, F->D) and gains a new trailing note in F.
$ws.Range("A14").Value = "Meetings:"
$ws.Range("C14").Value = "10.30 Init"
$ws.Range("D14").Value = "2:00 technical"
$ws.Range("F14").Value = "3-4 Mike"

# --- New "Meeting minutes" mini-table: bold headers in row 16, long wrapped
#     notes in row 17 (tall row). Values are written first (so the shared
#     strings land in the same order as the original edit), then formatting.
$ws.Range("C16").Value = "Meeting 1"
$ws.Range("D16").Value = "Meeting 2"
$ws.Range("C17").Value = "Discussed opening of project" + [char]10 + "Hours required" + [char]10 + "Slight changes to spec" + [char]10 + "What data is required "
$ws.Range("D17").Value = "Technical side of project" + [char]10 + "AngularJS for part1 with c# being involved as well" + [char]10 + "Met with other technical personal in the staff"

$ws.Range("C17").WrapText = $true
$ws.Range("C17").VerticalAlignment = -4160

$ws.Range("C16").Font.Bold = $true
$ws.Range("D16").Font.Bold = $true

$ws.Range("D17").WrapText = $true

$ws.Rows.Item(17).RowHeight = 130.5

# --- Back to the Week2 data row: fill in the previously-empty cells.
$ws.Range("D12").Value = "9->5"
$ws.Range("E12").Value = "9->12"
$ws.Range("E12").NumberFormat = "d-mmm"
$ws.Range("F12").Value = "9->12"

# --- New row 31: a trailing "Meetings:" label under the (still-empty) week-3
#     block that now lives at rows 27-30.
$ws.Range("A31").Value = "Meetings:"

# --- Column widths for the new notes columns / meeting table.
$ws.Columns.Item(3).ColumnWidth = 21.666666666666664
$ws.Columns.Item(4).ColumnWidth = 17.166666666666664
$ws.Columns.Item(12).ColumnWidth = 18.666666666666664

# --- Selection / view tidy-up to match the saved state.
$ws.Range("K14").Select() | Out-Null

# --- Page orientation.
$ws.PageSetup.Orientation = 1
